$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the login credential rows with new values
$ws.Range("A2").Value = "mngr212597"
$ws.Range("B2").Value = "urEguzu"
$ws.Range("A4").Value = "mngr212597"
$ws.Range("B4").Value = "urEguzu"

# Move active selection to B5
$ws.Range("B5").Select()
